$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 55
$ws.Range("B4").Value = 53
$ws.Range("B6").Value = 34
$ws.Range("B8").Value = 54
$ws.Range("B9").Value = 37
$ws.Range("B10").Value = 46
$ws.Range("B11").Value = 34
$ws.Range("B12").Value = 28
$ws.Range("C12").Value = 28
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 11
